# "Add files via upload" - re-upload of "Low-thoughput validation.xlsx"
#
# Net effect observed in the OOXML diff for the "Double mutants validation"
# sheet (sheet3.xml): the "Nscore" column (old column C) was removed, and the
# "Nham_aa" column (old column D) slid left to become the new column C.
# That is exactly what Excel does when you select the whole column C and
# delete it. Because the "Nscore" shared string is then unused anywhere in
# the workbook, Excel/the engine drops it from sharedStrings.xml and
# renumbers every other shared-string index down by one - which is why every
# other sheet that referenced a shared string index above "Nscore"'s slot
# shows a "-1" shift in the diff even though none of their actual values
# changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Double mutants validation")

# Delete the whole "Nscore" column; "Nham_aa" shifts left into column C.
$ws.Columns("C").Delete()

# Mirror the resulting selection/view state from the saved file: after
# deleting column C, Excel leaves the (now data-bearing) column C selected.
$ws.Activate()
$ws.Range("C1:C1048576").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
